# UP2P-261 - single file selection test. Multi file selection is not allowed.
# Adds a new "ImportFiledata" worksheet (between LogOutData and PrivacyPolicyData)
# that drives a single-file import/selection test, and nudges a couple of
# leftover cursor selections.

$wb = $excel.ActiveWorkbook

# --- Common: move the lingering selection cursor ----------------------------
$common = $wb.Worksheets.Item("Common")
$common.Activate()
$common.Range("B5").Select()

# --- SigninData: move the lingering selection cursor -------------------------
$signin = $wb.Worksheets.Item("SigninData")
$signin.Activate()
$signin.Range("A2").Select()

# --- New sheet: ImportFiledata, inserted right after LogOutData -------------
$afterSheet = $wb.Worksheets.Item("LogOutData")
$importSheet = $wb.Worksheets.Add($null, $afterSheet)
$importSheet.Name = "ImportFiledata"

# Header row
$importSheet.Range("A1").Value = "username"
$importSheet.Range("B1").Value = "password"
$importSheet.Range("C1").Value = "test_file_folder"
$importSheet.Range("D1").Value = "test_file_name"

# Data row - credentials pulled from Common, single test file target
$importSheet.Range("A2").Formula = "=Common!B7"
$importSheet.Range("B2").Formula = "=Common!B5"
$importSheet.Range("C2").Value = "TestData"
$importSheet.Range("D2").Value = "meals.csv"

# Column widths to fit the new headers
$importSheet.Columns.Item(1).ColumnWidth = 18.63
$importSheet.Columns.Item(2).ColumnWidth = 25.58
$importSheet.Columns.Item(3).ColumnWidth = 17.79
$importSheet.Columns.Item(4).ColumnWidth = 25.58

# Make this the active sheet/tab with its own cursor position, and leave it
# as the last-activated (tabSelected) sheet, matching the authored workbook.
$importSheet.Activate()
$importSheet.Range("B3").Select()
